# Daily attendance processing - 2025-12-23 19:07:57
# Normalize the "Recorded By" column (G) so that when the literal user
# "System" appears together with other recorders, it is listed last
# instead of first, e.g. "System, someone@example.com" -> "someone@example.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$prefix = "System, "

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val.ToString().StartsWith($prefix)) {
        $rest = $val.ToString().Substring($prefix.Length)
        $cell.Value = $rest + ", System"
    }
}
